# Update the "output generated at" snapshot: new date + new set of
# two-digit-by-one-digit division problems.

$d = $word.ActiveDocument

# --- Update the date/day heading (first paragraph) ---
$d.Content.Find.Execute("2025-10-25 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-10-26 Sunday", 2)

# --- Update the division problems inside the single table ---
# The table is 20 rows x 5 columns; only rows 1, 5, 9, 13, 17 hold the
# problem text (the other rows are blank answer rows). Addressing cells
# directly (rather than a global text replace) is required because some
# old values repeat (e.g. "41÷4=" appears twice in row 5).
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "16÷6=" },
    @{ Row = 1;  Col = 2; New = "57÷9=" },
    @{ Row = 1;  Col = 3; New = "64÷7=" },
    @{ Row = 1;  Col = 4; New = "79÷6=" },
    @{ Row = 1;  Col = 5; New = "29÷8=" },

    @{ Row = 5;  Col = 1; New = "80÷6=" },
    @{ Row = 5;  Col = 2; New = "33÷4=" },
    @{ Row = 5;  Col = 3; New = "57÷2=" },
    @{ Row = 5;  Col = 4; New = "68÷8=" },
    @{ Row = 5;  Col = 5; New = "32÷7=" },

    @{ Row = 9;  Col = 1; New = "65÷6=" },
    @{ Row = 9;  Col = 2; New = "14÷3=" },
    @{ Row = 9;  Col = 3; New = "54÷9=" },
    @{ Row = 9;  Col = 4; New = "74÷6=" },
    @{ Row = 9;  Col = 5; New = "61÷4=" },

    @{ Row = 13; Col = 1; New = "41÷8=" },
    @{ Row = 13; Col = 2; New = "79÷4=" },
    @{ Row = 13; Col = 3; New = "74÷7=" },
    @{ Row = 13; Col = 4; New = "40÷5=" },
    @{ Row = 13; Col = 5; New = "77÷9=" },

    @{ Row = 17; Col = 1; New = "86÷6=" },
    @{ Row = 17; Col = 2; New = "37÷7=" },
    @{ Row = 17; Col = 3; New = "83÷5=" },
    @{ Row = 17; Col = 4; New = "39÷5=" },
    @{ Row = 17; Col = 5; New = "83÷9=" }
)

foreach ($u in $updates) {
    $cell = $t.Rows.Item($u.Row).Cells.Item($u.Col)
    $cell.Range.Text = $u.New
}
